$d = $word.ActiveDocument

# Locate the existing first "Papers" bullet (Anderson B.M, Moore L., Bojechko C. ...)
# so we can insert a brand-new bullet immediately before it.
$targetRange = $d.Content
$targetRange.Find.Execute("Anderson B.M, Moore L., Bojechko C.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $targetRange.Paragraphs(1)
$start = $para.Range.Start

# Insert a new (empty) paragraph before the found one; it inherits the same
# paragraph formatting (ListParagraph style / bullet numbering) automatically.
$insertRange = $d.Range($start, $start)
$insertRange.InsertParagraphBefore()

# Cursor walks forward through the newly-created empty paragraph, one run at a time.
$pos = $start

# Run 1: "Anderson B.M" (bold)
$text = "Anderson B.M"
$cursor = $d.Range($pos, $pos)
$cursor.InsertAfter($text)
$cursor.SetRange($pos, $pos + $text.Length)
$cursor.Bold = 1
$cursor.BoldBi = 1
$pos = $pos + $text.Length

# Run 2: ", Bojechko C. "
$text = ", Bojechko C. "
$cursor = $d.Range($pos, $pos)
$cursor.InsertAfter($text)
$cursor.SetRange($pos, $pos + $text.Length)
$pos = $pos + $text.Length

# Run 3: title (italic)
$text = "DICOM Attribute Manipulation Tool: Easily Change Frame of Reference, Series Instance, and Study Instance UID"
$cursor = $d.Range($pos, $pos)
$cursor.InsertAfter($text)
$cursor.SetRange($pos, $pos + $text.Length)
$cursor.Italic = 1
$cursor.ItalicBi = 1
$pos = $pos + $text.Length

# Run 4: " "
$text = " "
$cursor = $d.Range($pos, $pos)
$cursor.InsertAfter($text)
$cursor.SetRange($pos, $pos + $text.Length)
$pos = $pos + $text.Length

# Run 5: "Practical Radiation Oncology (Submitted"
$text = "Practical Radiation Oncology (Submitted"
$cursor = $d.Range($pos, $pos)
$cursor.InsertAfter($text)
$cursor.SetRange($pos, $pos + $text.Length)
$pos = $pos + $text.Length

# Run 6: " 11/2024"
$text = " 11/2024"
$cursor = $d.Range($pos, $pos)
$cursor.InsertAfter($text)
$cursor.SetRange($pos, $pos + $text.Length)
$pos = $pos + $text.Length

# Run 7: ")"
$text = ")"
$cursor = $d.Range($pos, $pos)
$cursor.InsertAfter($text)
$cursor.SetRange($pos, $pos + $text.Length)
$pos = $pos + $text.Length
